# Update "Pais" (countries/provinces) COVID data sheet: refresh counts and re-sort ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title cell: timestamp of the data refresh.
$ws.Cells.Item(1,1).Value = "Datos actualizados a 7 de Abril de 2020 a las 23:52"

# Each entry: row number (already sorted by "Casos totales" desc), country name,
# then Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes.
$data = @(
  ,@(4, "Estados Unidos", 394182, 27178, 21571, 359895, 9169, 1845, 12716)
  ,@(5, "España", 140618, 3943, 43208, 83498, 7069, 571, 13912)
  ,@(6, "Italia", 135586, 3039, 24392, 94067, 3792, 604, 17127)
  ,@(7, "Francia", 109069, 11059, 19337, 79404, 7131, 1417, 10328)
  ,@(8, "Alemania", 107591, 4216, 36081, 69498, 4895, 202, 2012)
  ,@(9, "China", 81740, 32, 77167, 1242, 211, 0, 3331)
  ,@(10, "Iran", 62589, 2089, 27039, 31678, 3987, 133, 3872)
  ,@(11, "Reino Unido", 55242, 3634, 135, 48948, 1559, 786, 6159)
  ,@(12, "Turquia", 34109, 3892, 1582, 31802, 1474, 76, 725)
  ,@(13, "Suiza", 22253, 596, 8704, 12728, 391, 56, 821)
  ,@(14, "Belgica", 22194, 1380, 4157, 16002, 1260, 403, 2035)
  ,@(15, "Paises Bajos", 19580, 777, 250, 17229, 1424, 234, 2101)
  ,@(16, "Canada", 17872, 1205, 4028, 13467, 426, 54, 377)
  ,@(17, "Brasil", 13831, 1648, 127, 13023, 296, 117, 681)
  ,@(18, "Austria", 12639, 342, 4046, 8350, 243, 23, 243)
  ,@(19, "Portugal", 12442, 712, 184, 11913, 271, 34, 345)
  ,@(20, "Corea del Sur", 10331, 47, 6694, 3445, 55, 6, 192)
  ,@(21, "Israel", 9248, 344, 770, 8413, 149, 8, 65)
  ,@(22, "Suecia", 7693, 487, 205, 6897, 640, 114, 591)
  ,@(23, "Rusia", 7497, 1154, 494, 6945, 8, 11, 58)
  ,@(24, "Noruega", 5953, 88, 32, 5832, 78, 13, 89)
  ,@(25, "Australia", 5919, 24, 2547, 3324, 93, 3, 48)
  ,@(26, "Irlanda", 5709, 345, 25, 5474, 165, 36, 210)
  ,@(27, "India", 5311, 533, 421, 4740, 0, 14, 150)
  ,@(28, "Chile", 5116, 301, 898, 4175, 337, 6, 43)
  ,@(29, "Dinamarca", 5071, 390, 1491, 3377, 127, 16, 203)
  ,@(30, "Chequia", 4944, 122, 172, 4684, 86, 10, 88)
  ,@(31, "Polonia", 4848, 435, 191, 4528, 50, 22, 129)
  ,@(32, "Rumania", 4417, 360, 460, 3760, 274, 21, 197)
  ,@(33, "Pakistan", 4035, 269, 429, 3549, 28, 4, 57)
  ,@(34, "Malasia", 3963, 170, 1321, 2579, 92, 1, 63)
  ,@(35, "Japon", 3906, 0, 592, 3222, 79, 0, 92)
  ,@(36, "Filipinas", 3764, 104, 84, 3503, 1, 14, 177)
  ,@(37, "Ecuador", 3747, 0, 100, 3456, 156, 0, 191)
  ,@(38, "Luxemburgo", 2970, 127, 500, 2426, 35, 3, 44)
  ,@(39, "Peru", 2954, 393, 1301, 1546, 109, 15, 107)
  ,@(40, "Arabia Saudita", 2795, 190, 615, 2139, 41, 3, 41)
  ,@(41, "Indonesia", 2738, 247, 204, 2313, 0, 12, 221)
  ,@(42, "Serbia", 2447, 247, 118, 2268, 109, 3, 61)
  ,@(43, "Mexico", 2439, 296, 633, 1681, 89, 31, 125)
  ,@(44, "Emiratos Arabes Unidos", 2359, 283, 186, 2161, 1, 1, 12)
  ,@(45, "Finlandia", 2308, 132, 300, 1974, 81, 7, 34)
  ,@(46, "Tailandia", 2258, 38, 888, 1343, 61, 1, 27)
  ,@(47, "Panama", 2100, 0, 14, 2031, 88, 0, 55)
  ,@(48, "Catar", 2057, 225, 150, 1901, 37, 2, 6)
  ,@(49, "Republica Dominicana", 1956, 128, 36, 1822, 147, 12, 98)
  ,@(50, "Grecia", 1832, 77, 269, 1482, 90, 2, 81)
  ,@(51, "Colombia", 1780, 201, 100, 1630, 76, 4, 50)
  ,@(52, "Sudafrica", 1749, 63, 95, 1641, 7, 1, 13)
  ,@(53, "Argentina", 1628, 0, 338, 1234, 96, 3, 56)
  ,@(54, "Islandia", 1586, 24, 559, 1021, 11, 0, 6)
  ,@(55, "Singapur", 1481, 106, 377, 1098, 29, 0, 6)
  ,@(56, "Argelia", 1468, 45, 113, 1162, 46, 20, 193)
  ,@(57, "Ucrania", 1462, 143, 28, 1389, 16, 7, 45)
  ,@(58, "Egipto", 1450, 128, 276, 1080, 0, 9, 94)
  ,@(59, "Croacia", 1282, 60, 167, 1097, 35, 2, 18)
  ,@(60, "Marruecos", 1184, 64, 93, 1001, 1, 10, 90)
  ,@(61, "Nueva Zelanda", 1160, 54, 241, 918, 4, 0, 1)
  ,@(62, "Estonia", 1149, 41, 69, 1059, 12, 2, 21)
  ,@(63, "Irak", 1122, 91, 373, 684, 0, 1, 65)
  ,@(64, "Eslovenia", 1059, 38, 102, 921, 30, 6, 36)
  ,@(65, "Moldavia", 1056, 91, 40, 994, 80, 3, 22)
  ,@(66, "Hong Kong", 936, 21, 236, 696, 12, 0, 4)
  ,@(67, "Lituania", 880, 37, 8, 857, 11, 0, 15)
  ,@(68, "Bielorrusia", 861, 161, 54, 794, 31, 0, 13)
  ,@(69, "Armenia", 853, 20, 87, 758, 30, 0, 8)
  ,@(70, "Hungria", 817, 73, 71, 699, 17, 9, 47)
  ,@(71, "Barein", 811, 55, 458, 348, 4, 1, 5)
  ,@(72, "Bosnia y Herzegovina", 764, 90, 68, 663, 4, 4, 33)
  ,@(73, "Kuwait", 743, 78, 105, 637, 23, 0, 1)
  ,@(74, "Azerbaiyan", 717, 76, 44, 665, 23, 1, 8)
  ,@(75, "Crucero", 712, 0, 619, 82, 10, 0, 11)
  ,@(76, "Kazajistan", 697, 35, 51, 640, 16, 0, 6)
  ,@(77, "Camerun", 658, 0, 43, 606, 0, 0, 9)
  ,@(78, "Tunez", 623, 27, 25, 575, 39, 1, 23)
  ,@(79, "Republica de Macedonia", 599, 29, 30, 543, 15, 3, 26)
  ,@(80, "Eslovaquia", 581, 47, 13, 566, 3, 0, 2)
  ,@(81, "Bulgaria", 577, 28, 42, 512, 21, 1, 23)
  ,@(82, "Letonia", 548, 6, 16, 530, 5, 1, 2)
  ,@(83, "Libano", 548, 7, 62, 467, 29, 0, 19)
  ,@(84, "Principado de Andorra", 545, 20, 39, 484, 17, 1, 22)
  ,@(85, "Uzbekistan", 504, 47, 30, 472, 8, 0, 2)
  ,@(86, "Republica de Chipre", 494, 29, 47, 438, 13, 0, 9)
  ,@(87, "Costa Rica", 483, 16, 24, 457, 14, 0, 2)
  ,@(88, "Afganistan", 423, 56, 18, 391, 0, 3, 14)
  ,@(89, "Uruguay", 415, 0, 123, 286, 14, 0, 6)
  ,@(90, "Cuba", 396, 46, 27, 358, 15, 2, 11)
  ,@(91, "Burkina Faso", 384, 20, 127, 238, 0, 1, 19)
  ,@(92, "Albania", 383, 6, 131, 230, 7, 1, 22)
  ,@(93, "Taiwan", 376, 3, 61, 310, 0, 0, 5)
  ,@(94, "Oman", 371, 40, 67, 302, 3, 0, 2)
  ,@(95, "Jordania", 353, 4, 138, 209, 5, 0, 6)
  ,@(96, "Reunion", 349, 0, 40, 309, 4, 0, 0)
  ,@(97, "Costa de Marfil", 349, 26, 41, 305, 0, 0, 3)
  ,@(98, "Honduras", 305, 7, 6, 277, 10, 0, 22)
  ,@(99, "Malta", 293, 52, 5, 288, 4, 0, 0)
  ,@(100, "Ghana", 287, 73, 31, 251, 2, 0, 5)
  ,@(101, "Niger", 278, 25, 26, 241, 0, 1, 11)
  ,@(102, "San Marino", 277, 0, 35, 210, 14, 0, 32)
  ,@(103, "Mauricio", 268, 24, 8, 253, 3, 0, 7)
  ,@(104, "Estado de Palestina", 261, 7, 42, 218, 0, 0, 1)
  ,@(105, "Nigeria", 254, 16, 44, 204, 2, 1, 6)
  ,@(106, "Vietnam", 249, 4, 123, 126, 8, 0, 0)
  ,@(107, "Montenegro", 241, 8, 4, 235, 7, 0, 2)
  ,@(108, "Senegal", 237, 11, 105, 130, 1, 0, 2)
  ,@(109, "Kirguistan", 228, 12, 33, 191, 5, 0, 4)
  ,@(110, "Georgia", 196, 8, 46, 147, 6, 1, 3)
  ,@(111, "Bolivia", 194, 11, 2, 178, 3, 3, 14)
  ,@(112, "Sri Lanka", 185, 7, 42, 137, 5, 1, 6)
  ,@(113, "Islas Feroe", 184, 1, 120, 64, 2, 0, 0)
  ,@(114, "Consejo Danes para los Refugiados", 180, 19, 9, 153, 0, 0, 18)
  ,@(115, "Kenia", 172, 14, 7, 159, 2, 0, 6)
  ,@(116, "Venezuela", 165, 0, 65, 93, 6, 0, 7)
  ,@(117, "Mayotte", 164, 0, 15, 147, 3, 0, 2)
  ,@(118, "Banglades", 164, 41, 33, 114, 1, 5, 17)
  ,@(119, "Martinica", 152, 1, 50, 98, 20, 0, 4)
  ,@(120, "Isla de Man", 150, 11, 73, 76, 6, 0, 1)
  ,@(121, "Guinea", 144, 16, 5, 139, 0, 0, 0)
  ,@(122, "Guadalupe", 139, 0, 31, 101, 14, 0, 7)
  ,@(123, "Brunei", 135, 0, 85, 49, 3, 0, 1)
  ,@(124, "Paraguay", 115, 2, 15, 95, 1, 0, 5)
  ,@(125, "Camboya", 115, 1, 58, 57, 1, 0, 0)
  ,@(126, "Gibraltar", 113, 4, 60, 53, 0, 0, 0)
  ,@(127, "Trinidad yTobago", 107, 2, 1, 98, 0, 0, 8)
  ,@(128, "Ruanda", 105, 0, 7, 98, 0, 0, 0)
  ,@(129, "Republica de Yibuti", 90, 0, 9, 81, 0, 0, 0)
  ,@(130, "Madagascar", 88, 6, 2, 86, 6, 0, 0)
  ,@(131, "Monaco", 79, 2, 4, 74, 4, 0, 1)
  ,@(132, "El Salvador", 78, 9, 5, 69, 4, 0, 4)
  ,@(133, "Liechtenstein", 78, 1, 55, 22, 0, 0, 1)
  ,@(134, "Guatemala", 77, 7, 17, 57, 3, 0, 3)
  ,@(135, "Aruba", 74, 3, 14, 60, 0, 0, 0)
  ,@(136, "Guayana Francesa", 72, 0, 34, 38, 1, 0, 0)
  ,@(137, "Togo", 65, 7, 23, 39, 0, 0, 3)
  ,@(138, "Barbados", 63, 3, 6, 54, 4, 1, 3)
  ,@(139, "Jamaica", 59, 0, 8, 48, 0, 0, 3)
  ,@(140, "Mali", 56, 9, 12, 39, 0, 0, 5)
  ,@(141, "Uganda", 52, 0, 0, 52, 0, 0, 0)
  ,@(142, "Etiopia", 52, 8, 4, 46, 1, 0, 2)
  ,@(143, "Polinesia Francesa", 47, 5, 0, 47, 0, 0, 0)
  ,@(144, "Congo", 45, 0, 2, 38, 0, 0, 5)
  ,@(145, "Islas Caimanes", 45, 6, 6, 38, 0, 0, 1)
  ,@(146, "Macao", 44, 0, 10, 34, 1, 0, 0)
  ,@(147, "San Martin (Parte Holandesa)", 40, 3, 1, 33, 2, 0, 6)
  ,@(148, "Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
  ,@(149, "Zambia", 39, 0, 7, 31, 0, 0, 1)
  ,@(150, "Bermudas", 39, 0, 17, 20, 0, 0, 2)
  ,@(151, "Guinea-Bisau", 33, 15, 0, 33, 0, 0, 0)
  ,@(152, "Bahamas", 33, 0, 5, 22, 1, 1, 6)
  ,@(153, "Guyana", 33, 2, 8, 20, 8, 1, 5)
  ,@(154, "Guam", 32, 0, 0, 31, 0, 0, 1)
  ,@(155, "San Martin (Parte Francesa)", 32, 0, 7, 23, 6, 0, 2)
  ,@(156, "Eritrea", 31, 0, 0, 31, 0, 0, 0)
  ,@(157, "Gabon", 30, 6, 1, 28, 0, 0, 1)
  ,@(158, "Benin", 26, 0, 5, 20, 0, 0, 1)
  ,@(159, "Haiti", 25, 1, 0, 24, 0, 0, 1)
  ,@(160, "Tanzania", 24, 0, 5, 18, 0, 0, 1)
  ,@(161, "Birmania", 22, 0, 0, 21, 0, 0, 1)
  ,@(162, "Libia", 20, 1, 1, 18, 0, 0, 1)
  ,@(163, "Antigua y Barbuda", 19, 4, 0, 18, 1, 1, 1)
  ,@(164, "Siria", 19, 0, 3, 14, 0, 0, 2)
  ,@(165, "Maldivas", 19, 0, 13, 6, 0, 0, 0)
  ,@(166, "Nueva Caledonia", 18, 0, 1, 17, 0, 0, 0)
  ,@(167, "Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
  ,@(168, "Angola", 17, 1, 2, 13, 0, 0, 2)
  ,@(169, "Guinea Ecuatorial", 16, 0, 3, 13, 0, 0, 0)
  ,@(170, "Namibia", 16, 0, 3, 13, 0, 0, 0)
  ,@(171, "Fiyi", 15, 1, 0, 15, 0, 0, 0)
  ,@(172, "Dominica", 15, 0, 1, 14, 0, 0, 0)
  ,@(173, "Mongolia", 15, 0, 4, 11, 0, 0, 0)
  ,@(174, "Laos", 14, 2, 0, 14, 0, 0, 0)
  ,@(175, "Santa Lucia", 14, 0, 1, 13, 0, 0, 0)
  ,@(176, "Sudan", 14, 2, 2, 10, 0, 0, 2)
  ,@(177, "Liberia", 14, 0, 3, 8, 0, 0, 3)
  ,@(178, "Curazao", 13, 0, 7, 5, 0, 0, 1)
  ,@(179, "Granada", 12, 0, 0, 12, 2, 0, 0)
  ,@(180, "San Cristobal y Nieves", 11, 1, 0, 11, 0, 0, 0)
  ,@(181, "Seychelles", 11, 0, 0, 11, 0, 0, 0)
  ,@(182, "Zimbabue", 11, 1, 0, 9, 0, 1, 2)
  ,@(183, "Groenlandia", 11, 0, 10, 1, 0, 0, 0)
  ,@(184, "Mozambique", 10, 0, 1, 9, 0, 0, 0)
  ,@(185, "Surinam", 10, 0, 0, 9, 0, 0, 1)
  ,@(186, "Republica del Chad", 10, 1, 2, 8, 0, 0, 0)
  ,@(187, "Suazilandia", 10, 0, 4, 6, 0, 0, 0)
  ,@(188, "Nepal", 9, 0, 1, 8, 0, 0, 0)
  ,@(189, "Montserrat", 9, 0, 0, 7, 0, 0, 2)
  ,@(190, "Republica de Africa Central", 8, 0, 0, 8, 0, 0, 0)
  ,@(191, "Islas Turcas y Caicos", 8, 0, 0, 7, 0, 0, 1)
  ,@(192, "Somalia", 8, 1, 1, 7, 0, 0, 0)
  ,@(193, "Malaui", 8, 3, 0, 7, 1, 1, 1)
  ,@(194, "Santa Sede", 7, 0, 0, 7, 0, 0, 0)
  ,@(195, "Belice", 7, 0, 0, 6, 1, 0, 1)
  ,@(196, "San Vicente y las Granadinas", 7, 0, 1, 6, 0, 0, 0)
  ,@(197, "Cabo Verde", 7, 0, 1, 5, 0, 0, 1)
  ,@(198, "Sierra Leona", 6, 0, 0, 6, 0, 0, 0)
  ,@(199, "San Bartolome", 6, 0, 1, 5, 0, 0, 0)
  ,@(200, "Botsuana", 6, 0, 0, 5, 0, 0, 1)
  ,@(201, "Nicaragua", 6, 0, 0, 5, 0, 0, 1)
  ,@(202, "Mauritania", 6, 0, 2, 3, 0, 0, 1)
  ,@(203, "Butan", 5, 0, 2, 3, 0, 0, 0)
  ,@(204, "Santo Tome y Principe", 4, 0, 0, 4, 0, 0, 0)
  ,@(205, "Sahara Occidental", 4, 0, 0, 4, 0, 0, 0)
  ,@(206, "Gambia", 4, 0, 2, 1, 0, 0, 1)
  ,@(207, "Anguila", 3, 0, 0, 3, 0, 0, 0)
  ,@(208, "Islas Virgenes Britanicas", 3, 0, 0, 3, 0, 0, 0)
  ,@(209, "Burundi", 3, 0, 0, 3, 0, 0, 0)
  ,@(210, "Sudan del Sur", 2, 1, 0, 2, 0, 0, 0)
  ,@(211, "Islas Malvinas", 2, 0, 0, 2, 0, 0, 0)
  ,@(212, "Bonaire, San Eustaquio y Saba", 2, 0, 0, 2, 0, 0, 0)
  ,@(213, "Papua Nueva Guinea", 2, 0, 0, 2, 0, 0, 0)
  ,@(214, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
  ,@(215, "San Pedro y Miquelon", 1, 0, 0, 1, 0, 0, 0)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
    $ws.Cells.Item($r, 5).Value = $entry[5]
    $ws.Cells.Item($r, 6).Value = $entry[6]
    $ws.Cells.Item($r, 7).Value = $entry[7]
    $ws.Cells.Item($r, 8).Value = $entry[8]
}
